# This workbook (dlc12_iec61400-1ed3) lists DLC 1.2 load cases. Each case's
# wind direction offset lives in column F ([wdir]); column D is a cached
# formula string that embeds the (normalised, 0-360) wind direction inside
# the auto-generated case name, e.g.:
#   D = "dlc12_wsp"&TEXT(E,"00")&"_wdir"&TEXT(IF(F<0,F+360,F),"000")&"_s"&TEXT(G,"0000")
#
# The edit narrows the +/-10 degree yaw-misalignment cases to +/-8 degrees:
#   F = -10  ->  F = -8   (wdir 350 -> 352 once normalised to 0-360)
#   F =  10  ->  F =  8   (wdir 010 -> 008)
# Recalculating after the write refreshes the cached <v> in the D-column
# formula cells (and any other formulas, e.g. G/AN, that depend on F/E/G).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column F is the 6th column ([wdir]); find the last used row the same way
# Excel's Ctrl+Up would from the bottom of the sheet.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value()
    if ($val -eq -10) {
        $cell.Value = -8
    } elseif ($val -eq 10) {
        $cell.Value = 8
    }
}
